$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 171638
$ws.Cells.Item(2, 5).Value = 18251
$ws.Cells.Item(2, 6).Value = 18251
$ws.Cells.Item(2, 7).Value = 22538
$ws.Cells.Item(2, 8).Value = 17993
$ws.Cells.Item(2, 9).Value = 18012
$ws.Cells.Item(2, 10).Value = -19
$ws.Cells.Item(2, 11).Value = 279412
$ws.Cells.Item(2, 12).Value = 126930
$ws.Cells.Item(2, 13).Value = 152483
$ws.Cells.Item(2, 14).Value = 145067
$ws.Cells.Item(2, 15).Value = 7415
$ws.Cells.Item(2, 16).Value = 446
$ws.Cells.Item(2, 17).Value = 36774
$ws.Cells.Item(2, 18).Value = -36832
$ws.Cells.Item(2, 19).Value = -5594
$ws.Cells.Item(2, 20).Value = 30080
$ws.Cells.Item(2, 21).Value = 6694
$ws.Cells.Item(2, 22).Value = 69494
$ws.Cells.Item(2, 23).Value = 10.63
$ws.Cells.Item(2, 24).Value = 10.48
$ws.Cells.Item(2, 25).Value = 12.88
$ws.Cells.Item(2, 26).Value = 6.6
$ws.Cells.Item(2, 27).Value = 83.23999999999999
$ws.Cells.Item(2, 28).Value = 38317.34
$ws.Cells.Item(2, 29).Value = 22307
$ws.Cells.Item(2, 30).Value = 12.01
$ws.Cells.Item(2, 31).Value = 204504
$ws.Cells.Item(2, 32).Value = 1.31
$ws.Cells.Item(2, 33).Value = 9400
$ws.Cells.Item(2, 34).Value = 3.51
$ws.Cells.Item(2, 35).Value = 37.02
$ws.Cells.Item(2, 36).Value = 80745711

# Row 3
$ws.Cells.Item(3, 4).Value = 171367
$ws.Cells.Item(3, 5).Value = 17080
$ws.Cells.Item(3, 6).Value = 17080
$ws.Cells.Item(3, 7).Value = 20354
$ws.Cells.Item(3, 8).Value = 15159
$ws.Cells.Item(3, 9).Value = 15186
$ws.Cells.Item(3, 10).Value = -27
$ws.Cells.Item(3, 11).Value = 285814
$ws.Cells.Item(3, 12).Value = 132073
$ws.Cells.Item(3, 13).Value = 153741
$ws.Cells.Item(3, 14).Value = 152511
$ws.Cells.Item(3, 15).Value = 1230
$ws.Cells.Item(3, 16).Value = 446
$ws.Cells.Item(3, 17).Value = 37781
$ws.Cells.Item(3, 18).Value = -28805
$ws.Cells.Item(3, 19).Value = -9646
$ws.Cells.Item(3, 20).Value = 24788
$ws.Cells.Item(3, 21).Value = 12994
$ws.Cells.Item(3, 22).Value = 76440
$ws.Cells.Item(3, 23).Value = 9.970000000000001
$ws.Cells.Item(3, 24).Value = 8.85
$ws.Cells.Item(3, 25).Value = 10.21
$ws.Cells.Item(3, 26).Value = 5.36
$ws.Cells.Item(3, 27).Value = 85.91
$ws.Cells.Item(3, 28).Value = 40152.14
$ws.Cells.Item(3, 29).Value = 18807
$ws.Cells.Item(3, 30).Value = 11.46
$ws.Cells.Item(3, 31).Value = 215993
$ws.Cells.Item(3, 32).Value = 1
$ws.Cells.Item(3, 33).Value = 10000
$ws.Cells.Item(3, 34).Value = 4.64
$ws.Cells.Item(3, 35).Value = 46.63
$ws.Cells.Item(3, 36).Value = 80745711

# Row 4
$ws.Cells.Item(4, 4).Value = 170918
$ws.Cells.Item(4, 5).Value = 15357
$ws.Cells.Item(4, 6).Value = 15357
$ws.Cells.Item(4, 7).Value = 20961
$ws.Cells.Item(4, 8).Value = 16601
$ws.Cells.Item(4, 9).Value = 16760
$ws.Cells.Item(4, 10).Value = -159
$ws.Cells.Item(4, 11).Value = 312977
$ws.Cells.Item(4, 12).Value = 151812
$ws.Cells.Item(4, 13).Value = 161164
$ws.Cells.Item(4, 14).Value = 159714
$ws.Cells.Item(4, 15).Value = 1450
$ws.Cells.Item(4, 16).Value = 446
$ws.Cells.Item(4, 17).Value = 42432
$ws.Cells.Item(4, 18).Value = -24622
$ws.Cells.Item(4, 19).Value = -10448
$ws.Cells.Item(4, 20).Value = 24905
$ws.Cells.Item(4, 21).Value = 17527
$ws.Cells.Item(4, 22).Value = 76715
$ws.Cells.Item(4, 23).Value = 8.98
$ws.Cells.Item(4, 24).Value = 9.710000000000001
$ws.Cells.Item(4, 25).Value = 10.74
$ws.Cells.Item(4, 26).Value = 5.54
$ws.Cells.Item(4, 27).Value = 94.2
$ws.Cells.Item(4, 28).Value = 42270.33
$ws.Cells.Item(4, 29).Value = 20756
$ws.Cells.Item(4, 30).Value = 10.79
$ws.Cells.Item(4, 31).Value = 226194
$ws.Cells.Item(4, 32).Value = 0.99
$ws.Cells.Item(4, 33).Value = 10000
$ws.Cells.Item(4, 34).Value = 4.46
$ws.Cells.Item(4, 35).Value = 42.13
$ws.Cells.Item(4, 36).Value = 80745711

# Row 5
$ws.Cells.Item(5, 4).Value = 175200
$ws.Cells.Item(5, 5).Value = 15366
$ws.Cells.Item(5, 6).Value = 15366
$ws.Cells.Item(5, 7).Value = 34032
$ws.Cells.Item(5, 8).Value = 26576
$ws.Cells.Item(5, 9).Value = 25998
$ws.Cells.Item(5, 10).Value = 578
$ws.Cells.Item(5, 11).Value = 334287
$ws.Cells.Item(5, 12).Value = 153995
$ws.Cells.Item(5, 13).Value = 180292
$ws.Cells.Item(5, 14).Value = 178421
$ws.Cells.Item(5, 15).Value = 1871
$ws.Cells.Item(5, 16).Value = 446
$ws.Cells.Item(5, 17).Value = 38558
$ws.Cells.Item(5, 18).Value = -30706
$ws.Cells.Item(5, 19).Value = -8266
$ws.Cells.Item(5, 20).Value = 27159
$ws.Cells.Item(5, 21).Value = 11400
$ws.Cells.Item(5, 22).Value = 77717
$ws.Cells.Item(5, 23).Value = 8.77
$ws.Cells.Item(5, 24).Value = 15.17
$ws.Cells.Item(5, 25).Value = 15.38
$ws.Cells.Item(5, 26).Value = 8.210000000000001
$ws.Cells.Item(5, 27).Value = 85.41
$ws.Cells.Item(5, 28).Value = 46488.12
$ws.Cells.Item(5, 29).Value = 32198
$ws.Cells.Item(5, 30).Value = 8.289999999999999
$ws.Cells.Item(5, 31).Value = 252689
$ws.Cells.Item(5, 32).Value = 1.06
$ws.Cells.Item(5, 33).Value = 10000
$ws.Cells.Item(5, 34).Value = 3.75
$ws.Cells.Item(5, 35).Value = 27.16
$ws.Cells.Item(5, 36).Value = 80745711

# Row 6
$ws.Cells.Item(6, 4).Value = 168740
$ws.Cells.Item(6, 5).Value = 12018
$ws.Cells.Item(6, 6).Value = 12018
$ws.Cells.Item(6, 7).Value = 39760
$ws.Cells.Item(6, 8).Value = 31320
$ws.Cells.Item(6, 9).Value = 31279
$ws.Cells.Item(6, 11).Value = 423691
$ws.Cells.Item(6, 12).Value = 200199
$ws.Cells.Item(6, 13).Value = 223492
$ws.Cells.Item(6, 14).Value = 224708
$ws.Cells.Item(6, 16).Value = 446
$ws.Cells.Item(6, 17).Value = 43326
$ws.Cells.Item(6, 18).Value = -40477
$ws.Cells.Item(6, 19).Value = -2383
$ws.Cells.Item(6, 20).Value = 27924
$ws.Cells.Item(6, 21).Value = 15402
$ws.Cells.Item(6, 22).Value = 100761
$ws.Cells.Item(6, 23).Value = 7.12
$ws.Cells.Item(6, 24).Value = 18.56
$ws.Cells.Item(6, 25).Value = 15.52
$ws.Cells.Item(6, 26).Value = 8.26
$ws.Cells.Item(6, 27).Value = 89.58
$ws.Cells.Item(6, 28).Value = 56140.21
$ws.Cells.Item(6, 29).Value = 38738
$ws.Cells.Item(6, 30).Value = 6.96
$ws.Cells.Item(6, 31).Value = 312660
$ws.Cells.Item(6, 32).Value = 0.86
$ws.Cells.Item(6, 33).Value = 10000
$ws.Cells.Item(6, 34).Value = 3.71
$ws.Cells.Item(6, 35).Value = 22.94
$ws.Cells.Item(6, 36).Value = 80745711

# Row 7
$ws.Cells.Item(7, 4).Value = 179575
$ws.Cells.Item(7, 5).Value = 12197
$ws.Cells.Item(7, 7).Value = 14146
$ws.Cells.Item(7, 8).Value = 11022
$ws.Cells.Item(7, 9).Value = 11035
$ws.Cells.Item(7, 11).Value = 433104
$ws.Cells.Item(7, 12).Value = 205585
$ws.Cells.Item(7, 13).Value = 227520
$ws.Cells.Item(7, 14).Value = 228637
$ws.Cells.Item(7, 16).Value = 448
$ws.Cells.Item(7, 17).Value = 44949
$ws.Cells.Item(7, 18).Value = -35010
$ws.Cells.Item(7, 19).Value = -8379
$ws.Cells.Item(7, 20).Value = 35250
$ws.Cells.Item(7, 21).Value = 9471
$ws.Cells.Item(7, 23).Value = 6.79
$ws.Cells.Item(7, 24).Value = 6.14
$ws.Cells.Item(7, 25).Value = 4.87
$ws.Cells.Item(7, 26).Value = 2.57
$ws.Cells.Item(7, 27).Value = 90.36
$ws.Cells.Item(7, 29).Value = 13667
$ws.Cells.Item(7, 30).Value = 16.79
$ws.Cells.Item(7, 31).Value = 312617
$ws.Cells.Item(7, 32).Value = 0.73
$ws.Cells.Item(7, 33).Value = 10000
$ws.Cells.Item(7, 34).Value = 4.36
$ws.Cells.Item(7, 35).Value = 73.17

# Row 8
$ws.Cells.Item(8, 4).Value = 189325
$ws.Cells.Item(8, 5).Value = 13619
$ws.Cells.Item(8, 7).Value = 19617
$ws.Cells.Item(8, 8).Value = 15211
$ws.Cells.Item(8, 9).Value = 15304
$ws.Cells.Item(8, 11).Value = 442123
$ws.Cells.Item(8, 12).Value = 206794
$ws.Cells.Item(8, 13).Value = 235329
$ws.Cells.Item(8, 14).Value = 236650
$ws.Cells.Item(8, 16).Value = 448
$ws.Cells.Item(8, 17).Value = 45238
$ws.Cells.Item(8, 18).Value = -33426
$ws.Cells.Item(8, 19).Value = -9399
$ws.Cells.Item(8, 20).Value = 30162
$ws.Cells.Item(8, 21).Value = 15731
$ws.Cells.Item(8, 23).Value = 7.19
$ws.Cells.Item(8, 24).Value = 8.029999999999999
$ws.Cells.Item(8, 25).Value = 6.58
$ws.Cells.Item(8, 26).Value = 3.48
$ws.Cells.Item(8, 27).Value = 87.87
$ws.Cells.Item(8, 29).Value = 18954
$ws.Cells.Item(8, 30).Value = 12.11
$ws.Cells.Item(8, 31).Value = 323573
$ws.Cells.Item(8, 32).Value = 0.71
$ws.Cells.Item(8, 33).Value = 10081
$ws.Cells.Item(8, 34).Value = 4.39
$ws.Cells.Item(8, 35).Value = 53.19

# Row 9
$ws.Cells.Item(9, 4).Value = 198898
$ws.Cells.Item(9, 5).Value = 15169
$ws.Cells.Item(9, 7).Value = 26509
$ws.Cells.Item(9, 8).Value = 20630
$ws.Cells.Item(9, 9).Value = 20733
$ws.Cells.Item(9, 11).Value = 456313
$ws.Cells.Item(9, 12).Value = 208080
$ws.Cells.Item(9, 13).Value = 248233
$ws.Cells.Item(9, 14).Value = 249998
$ws.Cells.Item(9, 16).Value = 448
$ws.Cells.Item(9, 17).Value = 49988
$ws.Cells.Item(9, 18).Value = -37732
$ws.Cells.Item(9, 19).Value = -8105
$ws.Cells.Item(9, 20).Value = 31127
$ws.Cells.Item(9, 21).Value = 18308
$ws.Cells.Item(9, 23).Value = 7.63
$ws.Cells.Item(9, 24).Value = 10.37
$ws.Cells.Item(9, 25).Value = 8.52
$ws.Cells.Item(9, 26).Value = 4.59
$ws.Cells.Item(9, 27).Value = 83.81999999999999
$ws.Cells.Item(9, 29).Value = 25677
$ws.Cells.Item(9, 30).Value = 8.94
$ws.Cells.Item(9, 31).Value = 341824
$ws.Cells.Item(9, 32).Value = 0.67
$ws.Cells.Item(9, 33).Value = 10263
$ws.Cells.Item(9, 34).Value = 4.47
$ws.Cells.Item(9, 35).Value = 39.97
